$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that paired target "8.1" with indicators ['8.10.1','8.10.2'] (row 71)
# was mis-matched. Rows 71-78 shift up to take the values that used to be in
# rows 72-79 (targets 8.2 .. 8.9 with their matching indicators), and a new
# row 79 is created pairing the new target "8.10" with the indicator list
# that used to incorrectly sit on row 71 (['8.10.1','8.10.2']).

$ws.Range("B71").Value = "8.2"
$ws.Range("C71").Value = "['8.2.1']"

$ws.Range("B72").Value = "8.3"
$ws.Range("C72").Value = "['8.3.1']"

$ws.Range("B73").Value = "8.4"
$ws.Range("C73").Value = "['8.4.1', '8.4.2']"

$ws.Range("B74").Value = "8.5"
$ws.Range("C74").Value = "['8.5.1', '8.5.2']"

$ws.Range("B75").Value = "8.6"
$ws.Range("C75").Value = "['8.6.1']"

$ws.Range("B76").Value = "8.7"
$ws.Range("C76").Value = "['8.7.1']"

$ws.Range("B77").Value = "8.8"
$ws.Range("C77").Value = "['8.8.1', '8.8.2']"

$ws.Range("B78").Value = "8.9"
$ws.Range("C78").Value = "['8.9.1']"

$ws.Range("B79").Value = "8.10"
$ws.Range("C79").Value = "['8.10.1', '8.10.2']"

# Update the visible view to match: the new row 79 selected (full-row
# selection), scrolled so row 157 is near the top of the window.
$excel.ActiveWindow.ScrollRow = 157
$ws.Range("A79:XFD79").Select()
